$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.729.54"
$ws.Cells.Item(2, 5).Value = "  +0.59%  "

$ws.Cells.Item(3, 4).Value = "1.888.77"
$ws.Cells.Item(3, 5).Value = "  +0.25%  "

$ws.Cells.Item(4, 5).Value = "  +0.03%  "

$ws.Cells.Item(5, 4).Value = "'248.89"
$ws.Cells.Item(5, 5).Value = "  +0.25%  "

$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 5).Value = "  -0.01%  "

$ws.Cells.Item(7, 4).Value = "'0.4743"
$ws.Cells.Item(7, 5).Value = "  -0.53%  "

$ws.Cells.Item(8, 4).Value = "'0.2931"
$ws.Cells.Item(8, 5).Value = "  +0.27%  "

$ws.Cells.Item(9, 4).Value = "'0.06531"
$ws.Cells.Item(9, 5).Value = "  -0.09%  "

$ws.Cells.Item(10, 4).Value = "'21.97"
$ws.Cells.Item(10, 5).Value = "  -0.61%  "

$ws.Cells.Item(11, 4).Value = "'0.07797"
$ws.Cells.Item(11, 5).Value = "  +0.92%  "

$ws.Cells.Item(12, 4).Value = "'97.10"
$ws.Cells.Item(12, 5).Value = "  -0.88%  "

$ws.Cells.Item(13, 4).Value = "1.892.43"
$ws.Cells.Item(13, 5).Value = "  +0.42%  "

$ws.Cells.Item(14, 4).Value = "'0.7378"
$ws.Cells.Item(14, 5).Value = "  -0.51%  "

$ws.Cells.Item(15, 4).Value = "'5.251"
$ws.Cells.Item(15, 5).Value = "  +1.86%  "

$ws.Cells.Item(16, 4).Value = "'284.13"
$ws.Cells.Item(16, 5).Value = "  +3.47%  "

$ws.Cells.Item(17, 4).Value = "30.814.58"
$ws.Cells.Item(17, 5).Value = "  +0.88%  "

$ws.Cells.Item(18, 4).Value = "'13.19"
$ws.Cells.Item(18, 5).Value = "  -2.41%  "

$ws.Cells.Item(19, 4).Value = "'0.000007558"
$ws.Cells.Item(19, 5).Value = "  -0.23%  "

$ws.Cells.Item(20, 4).Value = "'1.000"
$ws.Cells.Item(20, 5).Value = "  +0.00%  "

$ws.Cells.Item(21, 4).Value = "2.142.02"
$ws.Cells.Item(21, 5).Value = "  +0.69%  "

$ws.Cells.Item(22, 4).Value = "'5.326"
$ws.Cells.Item(22, 5).Value = "  +1.17%  "

$ws.Cells.Item(23, 5).Value = "  +0.11%  "

$ws.Cells.Item(24, 4).Value = "'6.275"
$ws.Cells.Item(24, 5).Value = "  +1.09%  "

$ws.Cells.Item(25, 4).Value = "'9.233"
$ws.Cells.Item(25, 5).Value = "  -0.93%  "

$ws.Cells.Item(26, 4).Value = "'164.31"
$ws.Cells.Item(26, 5).Value = "  +0.52%  "

$ws.Cells.Item(27, 4).Value = "'18.94"
$ws.Cells.Item(27, 5).Value = "  -0.06%  "

$ws.Cells.Item(28, 4).Value = "'1.923"
$ws.Cells.Item(28, 5).Value = "  -1.27%  "

$ws.Cells.Item(29, 4).Value = "'1.342"
$ws.Cells.Item(29, 5).Value = "  -1.84%  "

$ws.Cells.Item(30, 4).Value = "'0.09733"
$ws.Cells.Item(30, 5).Value = "  -3.32%  "

$ws.Cells.Item(31, 4).Value = "'1.496"
$ws.Cells.Item(31, 5).Value = "  -1.42%  "

$ws.Cells.Item(32, 4).Value = "'4.295"
$ws.Cells.Item(32, 5).Value = "  -0.93%  "

$ws.Cells.Item(33, 4).Value = "'4.210"
$ws.Cells.Item(33, 5).Value = "  +2.20%  "

$ws.Cells.Item(34, 4).Value = "'0.04874"
$ws.Cells.Item(34, 5).Value = "  +0.99%  "

$ws.Cells.Item(35, 4).Value = "'1.128"
$ws.Cells.Item(35, 5).Value = "  -0.45%  "

$ws.Cells.Item(36, 4).Value = "'0.6991"
$ws.Cells.Item(36, 5).Value = "  -0.58%  "

$ws.Cells.Item(37, 4).Value = "'2.724"
$ws.Cells.Item(37, 5).Value = "  +0.43%  "

$ws.Cells.Item(38, 4).Value = "'0.01910"
$ws.Cells.Item(38, 5).Value = "  +2.09%  "

$ws.Cells.Item(39, 4).Value = "'2.812"
$ws.Cells.Item(39, 5).Value = "  +2.09%  "

$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(40, 4).Value = "'6.379"
$ws.Cells.Item(40, 5).Value = "  +0.72%  "

$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(41, 4).Value = "'76.07"
$ws.Cells.Item(41, 5).Value = "  +5.95%  "

$ws.Cells.Item(42, 4).Value = "'2.019"
$ws.Cells.Item(42, 5).Value = "  +0.98%  "

$ws.Cells.Item(43, 4).Value = "'0.4266"
$ws.Cells.Item(43, 5).Value = "  +0.87%  "

$ws.Cells.Item(44, 5).Value = "  -0.01%  "

$ws.Cells.Item(45, 4).Value = "'0.8366"
$ws.Cells.Item(45, 5).Value = "  -0.83%  "

$ws.Cells.Item(46, 4).Value = "'101.43"
$ws.Cells.Item(46, 5).Value = "  -1.47%  "

$ws.Cells.Item(47, 4).Value = "'9.496"
$ws.Cells.Item(47, 5).Value = "  +1.68%  "

$ws.Cells.Item(48, 4).Value = "'7.070"
$ws.Cells.Item(48, 5).Value = "  -0.56%  "

$ws.Cells.Item(49, 4).Value = "'35.70"
$ws.Cells.Item(49, 5).Value = "  +0.12%  "

$ws.Cells.Item(50, 4).Value = "'919.62"
$ws.Cells.Item(50, 5).Value = "  -0.03%  "

$ws.Cells.Item(51, 4).Value = "'0.05761"
$ws.Cells.Item(51, 5).Value = "  +1.94%  "
